$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C105").Value = 7295
$ws.Range("C106:C252").Value = 7293
